# Zambia model (OnSSET_cost_paramters.xlsx) update:
# Add mg_pv_hybrid / mg_wind_hybrid rows and an "average" row (with an
# AVERAGE formula) to the "for_NEST" sheet, widen its first column, and
# move the on-screen selections to where the author last left them.

$wb = $excel.ActiveWorkbook
$wsCost = $wb.Worksheets.Item("OnSSET_cost_paramters")
$wsNest = $wb.Worksheets.Item("for_NEST")

# --- new data rows on "for_NEST" ---------------------------------------
$wsNest.Range("A8").Value = "mg_pv_hybrid"
$wsNest.Range("B8").Value = 503

$wsNest.Range("A9").Value = "mg_wind_hybrid"
$wsNest.Range("B9").Value = 2800

$wsNest.Range("A10").Value = "average"
$wsNest.Range("B10").Formula = "=AVERAGE(B2:B9)"

# Column A is a little too narrow for "mg_wind_hybrid" - widen it.
$wsNest.Columns.Item(1).ColumnWidth = 12

# --- view/selection state, matching where the author left each sheet ---
# Leave the OnSSET_cost_paramters sheet scrolled/selected around row 39,
# then finish on for_NEST (which stays the active tab), selected at B11.
$wsCost.Activate()
$wsCost.Range("A39:B39").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$wsNest.Activate()
$wsNest.Range("B11").Select()
